$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "DATA_FIM_DT" in column M, reusing the same header style as the
# rest of row 1 (bold font + thin border + centered alignment).
$ws.Range("A1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "DATA_FIM_DT"

# DATA_FIM for row 2 was refreshed to the current run date.
$ws.Range("I2").Value = "13/11/25"

# Register the lowercase datetime format first (numFmtId 164) then drop it back to
# General so the slot isn't left attached to any cell - only the format code stays
# registered in the style table, unused.
$ws.Range("M2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("M2").NumberFormat = 0

# Now apply the real uppercase datetime format (numFmtId 165) used by the new
# DATA_FIM_DT column, with the underlying Excel date serial values.
$ws.Range("M2:M5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M2").Value = 45931
$ws.Range("M3").Value = 45870
$ws.Range("M4").Value = 45839
$ws.Range("M5").Value = 45974
